# Applies the F-column ('想去人数' / want-to-go count) updates captured
# in the commit 'Update gh-pages to output generated at 456a3b4'.
# Each of the four sheets (展览 / 演出 / 本地生活 / 全部类型) gets the
# same kind of small incremental bump to its F column values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1474
$ws.Range("F3").Value = 113
$ws.Range("F4").Value = 2113
$ws.Range("F5").Value = 7322
$ws.Range("F6").Value = 69
$ws.Range("F7").Value = 4732
$ws.Range("F8").Value = 6965
$ws.Range("F11").Value = 1470
$ws.Range("F12").Value = 850
$ws.Range("F13").Value = 155
$ws.Range("F14").Value = 42
$ws.Range("F17").Value = 150
$ws.Range("F20").Value = 1133
$ws.Range("F21").Value = 749
$ws.Range("F24").Value = 39
$ws.Range("F25").Value = 139
$ws.Range("F28").Value = 158
$ws.Range("F29").Value = 14
$ws.Range("F30").Value = 32
$ws.Range("F33").Value = 542
$ws.Range("F34").Value = 426
$ws.Range("F35").Value = 64
$ws.Range("F37").Value = 355
$ws.Range("F43").Value = 11

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F16").Value = 1731
$ws.Range("F20").Value = 8
$ws.Range("F21").Value = 200
$ws.Range("F23").Value = 138
$ws.Range("F26").Value = 631
$ws.Range("F31").Value = 845
$ws.Range("F32").Value = 982
$ws.Range("F33").Value = 602
$ws.Range("F43").Value = 71

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 110
$ws.Range("F5").Value = 853
$ws.Range("F6").Value = 661
$ws.Range("F7").Value = 287
$ws.Range("F8").Value = 1523
$ws.Range("F9").Value = 2403

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1474
$ws.Range("F4").Value = 853
$ws.Range("F5").Value = 113
$ws.Range("F6").Value = 661
$ws.Range("F7").Value = 661
$ws.Range("F8").Value = 7322
$ws.Range("F9").Value = 69
$ws.Range("F10").Value = 4732
$ws.Range("F11").Value = 6965
$ws.Range("F13").Value = 1470
$ws.Range("F15").Value = 850
$ws.Range("F16").Value = 155
$ws.Range("F17").Value = 1523
$ws.Range("F18").Value = 2403
$ws.Range("F19").Value = 200
$ws.Range("F20").Value = 42
$ws.Range("F22").Value = 138
$ws.Range("F23").Value = 150
$ws.Range("F25").Value = 1133
$ws.Range("F26").Value = 631
$ws.Range("F27").Value = 749
$ws.Range("F29").Value = 139
$ws.Range("F30").Value = 158
$ws.Range("F32").Value = 845
$ws.Range("F33").Value = 32
$ws.Range("F35").Value = 982
$ws.Range("F36").Value = 542
$ws.Range("F37").Value = 602
$ws.Range("F38").Value = 426
$ws.Range("F39").Value = 64
$ws.Range("F42").Value = 355
$ws.Range("F48").Value = 71

